$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4, pushing the old row 4 (10842374003) down to row 5.
# This is the "space to paste the parent codes" mentioned in the commit message.
$ws.Rows.Item(4).Insert()

# A leading apostrophe forces text storage so the long numeric-looking codes
# are kept as text (matching the workbook's original inline-string typing)
# instead of being reinterpreted as numbers.

# Row 2: code 10923690001 (was row 3's code), now flagged as an error lookup.
$ws.Range("A2").Value = "'10923690001"
$ws.Range("B2").Value = "Error: "
$ws.Range("C2").Value = "Precio no disponible"
$ws.Range("D2").Value = "Cantidad de imágenes no disponible"

# Row 3: code 10727609001 (was row 2's code), also flagged as an error lookup.
$ws.Range("A3").Value = "'10727609001"
$ws.Range("B3").Value = "Error: "
$ws.Range("C3").Value = "Precio no disponible"
$ws.Range("D3").Value = "Cantidad de imágenes no disponible"

# Row 4: brand-new parent code pasted into the freshly inserted row.
$ws.Range("A4").Value = "'10820236005"
$ws.Range("B4").Value = "Error: "
$ws.Range("C4").Value = "Precio no disponible"
$ws.Range("D4").Value = "Cantidad de imágenes no disponible"

# Row 5: the original row 4 code (10842374003), shifted down, also now an error lookup.
$ws.Range("A5").Value = "'10842374003"
$ws.Range("B5").Value = "Error: "
$ws.Range("C5").Value = "Precio no disponible"
$ws.Range("D5").Value = "Cantidad de imágenes no disponible"
